$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    '47+41=88',
    '28-0=28',
    '38+17=55',
    '47-44=3',
    '29+43=72',
    '16+25=41',
    '12-7=5',
    '98-3=95',
    '51+9=60',
    '41+2=43',
    '57-0=57',
    '30+39=69',
    '62+32=94',
    '28-12=16',
    '72-50=22',
    '45+49=94',
    '47+20=67',
    '50+32=82',
    '52+11=63',
    '91-84=7',
    '30-6=24',
    '7+67=74',
    '33+4=37',
    '54+42=96',
    '73-64=9',
    '16-10=6',
    '58-6=52',
    '97-85=12',
    '6+92=98',
    '2+71=73',
    '48+16=64',
    '75+4=79',
    '73-31=42',
    '41+7=48',
    '51+9=60',
    '61-58=3',
    '21+53=74',
    '11+2=13',
    '84-81=3',
    '45-41=4',
    '25+24=49',
    '55+32=87',
    '8+5=13',
    '78-30=48',
    '87-74=13',
    '42+52=94',
    '27+53=80',
    '31+67=98',
    '98-6=92',
    '76-39=37',
    '71-31=40',
    '58-21=37',
    '75-55=20',
    '36-6=30',
    '83-53=30',
    '33-4=29',
    '56-52=4',
    '24+32=56',
    '55-45=10',
    '32+36=68',
    '86+1=87',
    '27+68=95',
    '52-36=16',
    '2+67=69',
    '40-5=35',
    '71-67=4',
    '83-65=18',
    '54-3=51',
    '17+24=41',
    '31+59=90',
    '88-59=29',
    '58-31=27',
    '28+64=92',
    '67-43=24',
    '91-65=26',
    '86-48=38',
    '88-38=50',
    '65+29=94',
    '18+62=80',
    '45+46=91',
    '98-79=19',
    '49+9=58',
    '34+4=38',
    '98-55=43',
    '80-18=62',
    '71-11=60',
    '47+34=81',
    '65-13=52',
    '4+8=12',
    '18+46=64',
    '7+24=31',
    '17+59=76',
    '5+92=97',
    '49-31=18',
    '55+44=99',
    '81-79=2',
    '88-22=66',
    '85-57=28',
    '63-31=32',
    '78-65=13'
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"
